$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 54 (this shifts the former rows
# 54..172 down to 55..173, growing the used range to A1:R173).
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(54, 1).Value  = 11
$ws.Cells.Item(54, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(54, 3).Value  = "Bíobío"
$ws.Cells.Item(54, 4).Value  = 44533
$ws.Cells.Item(54, 5).Value  = 8
$ws.Cells.Item(54, 6).Value  = 100112017
$ws.Cells.Item(54, 7).Value  = "Apio"
$ws.Cells.Item(54, 8).Value  = "Americana (o)"
$ws.Cells.Item(54, 9).Value  = "Primera"
$ws.Cells.Item(54, 10).Value = 270
$ws.Cells.Item(54, 11).Value = 5500
$ws.Cells.Item(54, 12).Value = 6000
$ws.Cells.Item(54, 13).Value = 5722
$ws.Cells.Item(54, 14).Value = "`$/docena de matas"
$ws.Cells.Item(54, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(54, 16).Value = 954
$ws.Cells.Item(54, 17).Value = 6
$ws.Cells.Item(54, 18).Value = "Hortaliza"
